# iscatterplot.xlsx demo-data cleanup: rename samples/groups from the
# CtrlOE/PrrOE/RagKO/RagWT/siLuc/siPRR scheme to a generic WT/Mut scheme,
# and update the selection + column A width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D (group) / column A (sample) values, row by row (rows 2..19).
$groups = @("WT_1","WT_1","WT_1","Mut_1","Mut_1","Mut_1","WT_2","WT_2","WT_2","Mut_2","Mut_2","Mut_2","WT_3","WT_3","WT_3","Mut_3","Mut_3","Mut_3")
$samples = @("WT_1_Rep1","WT_1_Rep2","WT_1_Rep3","Mut_1_Rep1","Mut_1_Rep2","Mut_1_Rep3","WT_2_Rep1","WT_2_Rep2","WT_2_Rep3","Mut_2_Rep1","Mut_2_Rep2","Mut_2_Rep3","WT_3_Rep1","WT_3_Rep2","WT_3_Rep3","Mut_3_Rep1","Mut_3_Rep2","Mut_3_Rep3")

# Write column D (group) first for every row, then column A (sample) for
# every row -- matches the insertion order of the shared-strings table.
for ($i = 0; $i -lt $groups.Count; $i++) {
    $row = 2 + $i
    $ws.Range("D$row").Value = $groups[$i]
}
for ($i = 0; $i -lt $samples.Count; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $samples[$i]
}

# Column A got a touch narrower once the longer CtrlOE_REP#/siPRR_REP# style
# strings were replaced by the shorter WT_#_Rep#/Mut_#_Rep# labels.
$ws.Columns("A").ColumnWidth = 10.8333333333333

# Selection moved to A17:A19 (active cell A17).
$ws.Range("A17:A19").Select()
